# DOMA-7423: updated meter import example
#
# The example "Лист1" sheet lists sample meters for units at
# "г Москва, ул Тверская, д 1" (rows 2-7) and "...д 2" (rows 8-11).
# Rows 2-7 previously all used Unit type = "Flat". This update diversifies
# the Unit type column (C) for rows 3-6 so the example demonstrates every
# supported unit type, matching the set already used lower in the sheet:
#   C3: Flat -> Parking place
#   C4: Flat -> Apartment
#   C5: Flat -> Warehouse unit
#   C6: Flat -> Commercial unit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Parking place"
$ws.Range("C4").Value = "Apartment"
$ws.Range("C5").Value = "Warehouse unit"
$ws.Range("C6").Value = "Commercial unit"
